$d = $word.ActiveDocument

# --- Change 1: split "Metodo" sentence with a line break ---
$find1 = $d.Content.Find
$search1 = "Aulas teóricas e práticas, visitas técnicas e exercícios dirigidos. Avaliação baseada em provas, exercícios e trabalhos práticos e relatórios."
$replace1 = "Aulas teóricas e práticas, visitas técnicas e exercícios dirigidos. ^lAvaliação baseada em provas, exercícios e trabalhos práticos e relatórios."
$result1 = $find1.Execute($search1, $true, $false, $false, $false, $false, $true, 1, $false, $replace1, 2)
Write-Output "Metodo replace result: $result1"
if (-not $result1) { throw "Failed to apply Metodo line-break split" }

# --- Change 2: split bibliography block into one entry per line ---
$find2 = $d.Content.Find
$search2 = "Bibliografia básica:AGRA FILHO, S,S. Planejamento e Gestão Ambiental no Brasil. Os Instrumentos da Política Nacional do Meio Ambiente, Rio de Janeiro, Elsevier, 2014FRANCO, M.A.R., Planejamento ambiental para a cidade sustentável, Ed. Annablume, 2000DEAK, C., SHIFFER, S.T.R., O processo de urbanização no Brasil, EDUSP, 1999IBGE, Instituto Brasileiro de Geografia e Estatística. Indicadores de Desenvolvimento Sustentável. Rio de Janeiro, IBGE, 2012.MOTA, S., Urbanização e meio ambiente, ABES Associação Brasileira de Engenharia Sanitária, 1999MENEZES, C.L., Desenvolvimento urbano e meio ambiente, Papirus, 1996PHILLIPI, Jr.A; MALHEIROS, T.F. Indicadores de Sustentabilidade e Gestão Ambiental. Editora Manole, 2012.SANTOS, M. A Urbanização Brasileira. 3 ed. São Paulo: HUCITEC, 1993. 155pSANTOS, R.F., Planejamento ambiental: teoria e prática, Editora Oficina de textos, 2004SECCHI, L. Análise de Políticas Públicas. Diagnóstico de Problemas, Recomendações de Soluções., São Paulo, Cengage Learning, 2016SOUZA, M.L. Mudar a Cidade: Uma introdução crítica ao planejamento e à gestão urbanos. Rio de Janeiro, Bertrand Brasil, 2003.VILLAÇA, F. Uma contribuição para a história do planejamento urbano no Brasil. In: DEAK, C; SCHIFFER, S.R (org) O processo de urbanização no Brasil. São Paulo, EDUSP, 1999.Bibliografia complementar:ALLEN, A., YOU, N., Sustainable urbanization – bridging the green and brown agendas, DPU, University College London, 2002ACSELRAD, H., Conflitos ambientais no Brasil, Fundação Henrich Boll, 2004BARDET, G., O urbanismo, Papirus, 1990BUARQUE, S.C., LIMA, R.R.A.; Manual de estratégia de desenvolvimento para aglomerações urbanas, Brasília, IPEA, 2005MENEGAT, R; ALMEIDA, G. Desenvolvimento Sustentável e Gestão Ambiental nas Cidades. Porto Alegre, Editora UFRGS, 2004."
$replace2 = "Bibliografia básica:^lAGRA FILHO, S,S. Planejamento e Gestão Ambiental no Brasil. Os Instrumentos da Política Nacional do Meio Ambiente, Rio de Janeiro, Elsevier, 2014^lFRANCO, M.A.R., Planejamento ambiental para a cidade sustentável, Ed. Annablume, 2000^lDEAK, C., SHIFFER, S.T.R., O processo de urbanização no Brasil, EDUSP, 1999^lIBGE, Instituto Brasileiro de Geografia e Estatística. Indicadores de Desenvolvimento Sustentável. Rio de Janeiro, IBGE, 2012.^lMOTA, S., Urbanização e meio ambiente, ABES Associação Brasileira de Engenharia Sanitária, 1999^lMENEZES, C.L., Desenvolvimento urbano e meio ambiente, Papirus, 1996^lPHILLIPI, Jr.A; MALHEIROS, T.F. Indicadores de Sustentabilidade e Gestão Ambiental. Editora Manole, 2012.^lSANTOS, M. A Urbanização Brasileira. 3 ed. São Paulo: HUCITEC, 1993. 155p^lSANTOS, R.F., Planejamento ambiental: teoria e prática, Editora Oficina de textos, 2004^lSECCHI, L. Análise de Políticas Públicas. Diagnóstico de Problemas, Recomendações de Soluções., São Paulo, Cengage Learning, 2016^lSOUZA, M.L. Mudar a Cidade: Uma introdução crítica ao planejamento e à gestão urbanos. Rio de Janeiro, Bertrand Brasil, 2003.^lVILLAÇA, F. Uma contribuição para a história do planejamento urbano no Brasil. In: DEAK, C; SCHIFFER, S.R (org) O processo de urbanização no Brasil. São Paulo, EDUSP, 1999.^l^lBibliografia complementar:^lALLEN, A., YOU, N., Sustainable urbanization – bridging the green and brown agendas, DPU, University College London, 2002^lACSELRAD, H., Conflitos ambientais no Brasil, Fundação Henrich Boll, 2004^lBARDET, G., O urbanismo, Papirus, 1990^lBUARQUE, S.C., LIMA, R.R.A.; Manual de estratégia de desenvolvimento para aglomerações urbanas, Brasília, IPEA, 2005^lMENEGAT, R; ALMEIDA, G. Desenvolvimento Sustentável e Gestão Ambiental nas Cidades. Porto Alegre, Editora UFRGS, 2004."
$result2 = $find2.Execute($search2, $true, $false, $false, $false, $false, $true, 1, $false, $replace2, 2)
Write-Output "Bibliografia replace result: $result2"
if (-not $result2) { throw "Failed to apply Bibliografia line-break split" }
